$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''60.884.80'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.22%  '
$ws.Range("D3").Value = '''2.917.78'
$ws.Range("D3").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''583.84'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.48%  '
$ws.Range("D6").Value = '''144.26'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.06%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("E8").Value = '  -2.54%  '
$ws.Range("D9").Value = '''2.916.31'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.61%  '
$ws.Range("D10").Value = '''6.85'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +7.07%  '
$ws.Range("E11").Value = '  -3.88%  '
$ws.Range("E12").Value = '  -4.00%  '
$ws.Range("E13").Value = '  -3.22%  '
$ws.Range("D14").Value = '''33.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -5.20%  '
$ws.Range("E15").Value = '  +0.27%  '
$ws.Range("D16").Value = '''3.400.64'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.83%  '
$ws.Range("D17").Value = '''60.828.60'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -3.21%  '
$ws.Range("D18").Value = '''6.74'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -4.99%  '
$ws.Range("D19").Value = '''2.917.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -3.75%  '
$ws.Range("D20").Value = '''431.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.39%  '
$ws.Range("D21").Value = '''13.64'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.47%  '
$ws.Range("D22").Value = '''0.683'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.72%  '
$ws.Range("D23").Value = '''7.14'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.50%  '
$ws.Range("D24").Value = '''80.48'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.14%  '
$ws.Range("E25").Value = '  -1.24%  '
$ws.Range("D26").Value = '''2.21'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.03%  '
$ws.Range("D27").Value = '''11.93'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -3.20%  '
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("D30").Value = '''7.26'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.50%  '
$ws.Range("E31").Value = '  -2.96%  '
$ws.Range("D32").Value = '''2.17'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.47%  '
$ws.Range("D33").Value = '''26.55'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.63%  '
$ws.Range("D34").Value = '''0.107'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.71%  '
$ws.Range("D35").Value = '''0.0₃0872'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.95%  '
$ws.Range("E36").Value = '  -2.55%  '
$ws.Range("D37").Value = '''5.66'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -4.35%  '
$ws.Range("D38").Value = '''3.02'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.74%  '
$ws.Range("D39").Value = '''49.76'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -1.54%  '
$ws.Range("E40").Value = '  -1.80%  '
$ws.Range("D41").Value = '''2.01'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -4.44%  '
$ws.Range("D42").Value = '''8.70'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.26%  '
$ws.Range("D43").Value = '''0.296'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -3.69%  '
$ws.Range("D44").Value = '''41.27'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.23%  '
$ws.Range("D45").Value = '''0.0349'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.65%  '
$ws.Range("D46").Value = '''376.58'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -4.62%  '
$ws.Range("D47").Value = '''2.674.14'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.14%  '
$ws.Range("D48").Value = '''132.31'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.27%  '
$ws.Range("D50").Value = '''24.57'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.09%  '
